$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 86) to the daily log table.
# Column A holds a date-like text value ("2025/10/10"); prefix it with an
# apostrophe so the engine stores it as literal text (matching the other
# rows in the sheet, which store dates as plain strings) instead of
# auto-converting it into a date serial number. Then strip the resulting
# "quote prefix" cell formatting so the cell ends up with the sheet's
# default (unstyled) appearance, just like the rest of the table.
$ws.Range("A86").Value = "'2025/10/10"
$ws.Range("A86").ClearFormats()

$ws.Range("B86").Value = "金"
$ws.Range("C86").Value = 1
$ws.Range("D86").Value = 29
